# Applies the "market page created and Top Coins moved, small changes to
# login and navbar visibility" commit to the Tuntikirjanpito workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# ---------------------------------------------------------------------
# New diary entries (rows 76-79), appended after the last existing entry
# (row 75). Column A only gets a date on the first row of a given day,
# column B is hours spent, column C is the description, column D is the
# project/category tag ("client").
# ---------------------------------------------------------------------

$ws.Cells.Item(76, 1).Value = 44578
$ws.Range("A75").Copy()
$ws.Range("A76").PasteSpecial(-4122)
$ws.Cells.Item(76, 2).Value = 1
$ws.Cells.Item(76, 3).Value = "NavBar profile näkyy vain kun käyttäjä kirjautuneena, TopCoins siirretty ja uudelleennimetty Market sivulle"
$ws.Cells.Item(76, 4).Value = "client"

$ws.Cells.Item(77, 2).Value = 1
$ws.Cells.Item(77, 3).Value = "uusi Gecko API linkki haettu ja testattu, top 5 kolikot market cap, hinta per coin, volume"
$ws.Cells.Item(77, 4).Value = "client"

$ws.Cells.Item(78, 2).Value = 1
$ws.Cells.Item(78, 3).Value = "Market sivun komponettien pohjien rakentelua, CoinCard, CoinCardList, NoData, MarketContainer"
$ws.Cells.Item(78, 4).Value = "client"

$ws.Cells.Item(79, 2).Value = 1
$ws.Cells.Item(79, 3).Value = "CoinCard tuunausta, Font Awsome -kirjaston lisäys"
$ws.Cells.Item(79, 4).Value = "client"

# ---------------------------------------------------------------------
# Update the totals row: formula now sums through the new last entry
# row (79) instead of 75.
# ---------------------------------------------------------------------

$ws.Range("B86").Formula = "=SUM(B2:B79)"

# ---------------------------------------------------------------------
# Restore the active cell selection left behind when the file was saved.
# ---------------------------------------------------------------------

$ws.Range("B80").Select()
